$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update "current as of" day label (merged R1:R2) from D5.2 to D6.1
$ws.Range("R1").Value = "D6.1"

# Row 15 (GBU-12): record 2 units expended on day D6.1 (column N)
$ws.Range("N15").Value = 2

# Row 18 (AGM-65 all types): record 16 units expended on day D6.1 (column N)
$ws.Range("N18").Value = 16

# Row 18's inventory total now drops into the "critical" range, so match
# the highlight formatting already used by the other low-stock rows
# (e.g. row 15) by copying formats from R15 onto R18.
$ws.Range("R15").Copy() | Out-Null
$ws.Range("R18").PasteSpecial(-4122) | Out-Null

# Leave the selection where the author left it when saving
$ws.Range("U8").Select() | Out-Null
